$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AB2").Value = 6.2
$ws.Range("AC2").Value = 11.5
$ws.Range("AF2").Value = 6.8
$ws.Range("AI2").Value = 260
$ws.Range("AJ2").Value = 11.5
$ws.Range("AM2").Value = 390
$ws.Range("F2").Value = 1.4
$ws.Range("T2").Value = 2.56
$ws.Range("AH3").Value = 17.5
$ws.Range("AJ3").Value = 46
$ws.Range("AK3").Value = 34
$ws.Range("AO3").Value = 27
$ws.Range("G3").Value = 2.92
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 3.45
$ws.Range("O3").Value = 1.34
$ws.Range("Q3").Value = 2.02
$ws.Range("T3").Value = 1.79
$ws.Range("U3").Value = 2.18
$ws.Range("X3").Value = 14
$ws.Range("Z3").Value = 18
$ws.Range("AB4").Value = 13.5
$ws.Range("AK4").Value = 13
$ws.Range("AN4").Value = 4.5
$ws.Range("H4").Value = 8.4
$ws.Range("I4").Value = 8.6
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 7
$ws.Range("R4").Value = 1.82
$ws.Range("U4").Value = 2.4
$ws.Range("X4").Value = 36
$ws.Range("F5").Value = 2.2
$ws.Range("I5").Value = 3.55
$ws.Range("K5").Value = 3.85
$ws.Range("N5").Value = 4.7
$ws.Range("P5").Value = 2.22
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 1.48
$ws.Range("S5").Value = 2.98
$ws.Range("Y5").Value = 16.5
$ws.Range("AA6").Value = 180
$ws.Range("AB6").Value = 7.6
$ws.Range("AE6").Value = 1000
$ws.Range("AJ6").Value = 20
$ws.Range("AM6").Value = 150
$ws.Range("O6").Value = 1.4
$ws.Range("S6").Value = 4.1
$ws.Range("U6").Value = 1.91
$ws.Range("X6").Value = 12
$ws.Range("AJ7").Value = 65
$ws.Range("H7").Value = 2.4
$ws.Range("P7").Value = 1.96
$ws.Range("U7").Value = 2.22
$ws.Range("AB8").Value = 8.4
$ws.Range("AC8").Value = 8
$ws.Range("AF8").Value = 11.5
$ws.Range("AN8").Value = 14.5
$ws.Range("AO8").Value = 110
$ws.Range("G8").Value = 1.91
$ws.Range("Q8").Value = 2.08
$ws.Range("Z8").Value = 36
$ws.Range("P9").Value = 2.12
$ws.Range("T9").Value = 1.72
$ws.Range("U9").Value = 2.32
$ws.Range("Y9").Value = 11.5
$ws.Range("AA10").Value = 70
$ws.Range("F10").Value = 2.46
$ws.Range("G10").Value = 2.48
$ws.Range("H10").Value = 3.45
$ws.Range("O10").Value = 1.44
$ws.Range("U10").Value = 1.95
$ws.Range("AN11").Value = 27
$ws.Range("F11").Value = 2.32
$ws.Range("G11").Value = 2.38
$ws.Range("I11").Value = 3.85
$ws.Range("Z11").Value = 26
$ws.Range("AA12").Value = 12
$ws.Range("AF12").Value = 85
$ws.Range("AG12").Value = 34
$ws.Range("AO12").Value = 5.7
$ws.Range("G12").Value = 9.6
$ws.Range("H12").Value = 1.42
$ws.Range("I12").Value = 1.43
$ws.Range("K12").Value = 5.4
$ws.Range("Q12").Value = 1.67
$ws.Range("R12").Value = 1.55
$ws.Range("T12").Value = 1.95
$ws.Range("AC13").Value = 8.800000000000001
$ws.Range("AF13").Value = 11
$ws.Range("AG13").Value = 9.800000000000001
$ws.Range("AH13").Value = 19.5
$ws.Range("AI13").Value = 75
$ws.Range("AJ13").Value = 17.5
$ws.Range("AL13").Value = 34
$ws.Range("AM13").Value = 110
$ws.Range("AN13").Value = 9.6
$ws.Range("F13").Value = 1.73
$ws.Range("G13").Value = 1.75
$ws.Range("I13").Value = 5.8
$ws.Range("S13").Value = 3.2
$ws.Range("X13").Value = 16
$ws.Range("Z13").Value = 55
$ws.Range("AG14").Value = 21
$ws.Range("AJ14").Value = 140
$ws.Range("AN14").Value = 65
$ws.Range("H14").Value = 1.77
$ws.Range("Y14").Value = 9.6
